$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28, shifting the existing
# rows 28-35 down to 29-36 (this also grows the sheet dimension to R36
# and copies formatting, e.g. the date style on column D, from the row
# above into the new row).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value2 = 44825
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112026
$ws.Range("G28").Value = "Haba"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14500
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 580
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
